$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 264: duplicate of the existing "Brant Hurter" entry (row 261) ---
$ws.Range("A264").Value = "Brant Hurter"
$ws.Range("B264").Value = "Brant Hurter"
$ws.Range("C264").Value = "https://midfield.mlbstatic.com/v1/people/676428/spots/120"
$ws.Range("D264").Value = "LHP"

# --- Rows 265-266: two brand-new pitchers ---
# Set the name cells first (for both new rows) so the new shared strings are
# interned in the same order the source workbook used: Ky Bush, Davis Martin,
# then their two headshot URLs.
$ws.Range("A265").Value = "Ky Bush"
$ws.Range("B265").Value = "Ky Bush"
$ws.Range("A266").Value = "Davis Martin"
$ws.Range("B266").Value = "Davis Martin"

$ws.Range("C265").Value = "https://midfield.mlbstatic.com/v1/people/681066/spots/120"
$ws.Range("C266").Value = "https://midfield.mlbstatic.com/v1/people/663436/spots/120"

$ws.Range("D265").Value = "LHP"
$ws.Range("D266").Value = "RHP"

# --- Hyperlinks on column C, matching the existing rows' style ---
$ws.Hyperlinks.Add($ws.Range("C264"), "https://midfield.mlbstatic.com/v1/people/676428/spots/120")
$ws.Hyperlinks.Add($ws.Range("C265"), "https://midfield.mlbstatic.com/v1/people/681066/spots/120")
$ws.Hyperlinks.Add($ws.Range("C266"), "https://midfield.mlbstatic.com/v1/people/663436/spots/120")

# Re-apply the sheet's normal Hyperlink cell style (Hyperlinks.Add nudges in
# its own variant xf) so the new cells land back on the same style index the
# rest of column C already uses.
$ws.Range("C264").Style = "Hyperlink"
$ws.Range("C265").Style = "Hyperlink"
$ws.Range("C266").Style = "Hyperlink"

# Match the selection recorded in the saved workbook.
$ws.Range("B264:B266").Select()
